$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 2 for the new IRCON entry. ---
# This shifts the existing row 2 (BLS) down to row 3, and existing
# row 3 (CRAFTSMAN) down to row 4 -- matching the target layout.
$ws.Rows.Item(2).Insert()

# --- New header cell M1 (Profit/Loss, same column meaning as G1) ---
$ws.Range("M1").Value = "Profit/Loss"

# --- Row 2 (new): 29-01-2024 / IRCON trade ---
$ws.Range("A2").Value = "29-01-2024"
$ws.Range("B2").Value = "IRCON"
$ws.Range("C2").Value = 247.25
$ws.Range("D2").Value = 238.8
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 20
$ws.Range("K2").Value = "Delivery"
$ws.Range("L2").Value = 0

# Establish the date number format (mm-dd-yy) once on A2, then reuse
# that same style for the other date cells below via copy/paste-special
# so the workbook ends up with a single shared style record instead of
# one per cell.
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy()

# --- Row 5 (new): BLS delivery, no exit yet ---
$ws.Range("A5").Value = 45293
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("B5").Value = "BLS"
$ws.Range("C5").Value = 420.75
$ws.Range("E5").Value = 40
$ws.Range("K5").Value = "Delivery"
$ws.Range("L5").Value = 0

# --- Row 6 (new): ELECTCAST delivery ---
$ws.Range("A6").Value = 45293
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").Value = "ELECTCAST"
$ws.Range("C6").Value = 171.35
$ws.Range("D6").Value = 177
$ws.Range("E6").Value = 19
$ws.Range("F6").Value = 19
$ws.Range("K6").Value = "Delivery"
$ws.Range("L6").Value = 0

# --- Row 7 (new): NHPC delivery, ends in a loss ---
$ws.Range("A7").Value = 45293
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B7").Value = "NHPC"
$ws.Range("C7").Value = 90.5
$ws.Range("D7").Value = 97
$ws.Range("E7").Value = 40
$ws.Range("F7").Value = 40
$ws.Range("K7").Value = "Delivery"
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = "Loss"

# --- Row 8 (new): CRAFTSMAN intraday ---
$ws.Range("A8").Value = 45324
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").Value = "CRAFTSMAN"
$ws.Range("C8").Value = 4291.2
$ws.Range("D8").Value = 4312.75
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = -21.55
$ws.Range("K8").Value = "Intraday"
$ws.Range("L8").Value = 0

# --- Row 9 (new): SUZLON delivery, symbol entered with a leading quote ---
$ws.Range("A9").Value = 45324
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B9").Value = "'SUZLON"
$ws.Range("C9").Value = 48.9
$ws.Range("E9").Value = 90
$ws.Range("K9").Value = "Delivery"
$ws.Range("L9").Value = 0

$ws.Application.CutCopyMode = 0

# --- Selection matches where the user was working next (row 10) ---
$ws.Range("G10").Select()
